$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": update Blancos/Reprobados/Aprobados/Por_Apro
#     and add Promedio (column H) for rows 2-4 ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2
$ws2P.Range("D2").Value = 4
$ws2P.Range("E2").Value = 4
$ws2P.Range("F2").Value = 21
$ws2P.Range("G2").Value = 84
$ws2P.Range("H2").Value = 7.5

# Row 3
$ws2P.Range("D3").Value = 7
$ws2P.Range("E3").Value = 7
$ws2P.Range("F3").Value = 14
$ws2P.Range("G3").Value = 66.67
$ws2P.Range("H3").Value = 8.1

# Row 4
$ws2P.Range("D4").Value = 4
$ws2P.Range("E4").Value = 4
$ws2P.Range("F4").Value = 24
$ws2P.Range("G4").Value = 85.70999999999999
$ws2P.Range("H4").Value = 8

# --- Sheet "Estadisticos Final": update Promedio (column H) for rows 2-3 ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("H2").Value = 7.5
$wsFinal.Range("H3").Value = 6.8
